$wb = $excel.ActiveWorkbook

# Sheets affected: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

# F-column (想去人数) updates for rows 3..17 (row 2 handled separately for column G)
$fUpdates = @{
    3  = 518
    4  = 16
    5  = 98
    8  = 58
    9  = 481
    10 = 6303
    11 = 163
    12 = 118
    13 = 1011
    14 = 302
    15 = 84
    16 = 173
    17 = 427
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # G2: number 50 -> text "不可售"
    $ws.Range("G2").Value = "不可售"

    # F column updates
    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }
}
